$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.181.61'
$ws.Range('E2').Value = '  +0.60%  '
$ws.Range('D3').Value = '1.842.38'
$ws.Range('E3').Value = '  +0.56%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6271'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.002'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07562'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2953'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.13%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.42'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07717'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('D12').Value = '1.841.72'
$ws.Range('E12').Value = '  +0.47%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.039'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.85%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6800'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.37%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '83.36'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.000009325'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.994'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.99%  '
$ws.Range('D18').Value = '29.166.00'
$ws.Range('E18').Value = '  +0.42%  '
$ws.Range('D19').Value = '2.088.49'
$ws.Range('E19').Value = '  +0.26%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '232.61'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.70%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.75'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.002'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.196'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.002'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '160.92'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.40%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1411'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.42%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.571'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.43%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.01'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.500'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.18%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.202'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.93%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.163'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.65%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05588'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.59%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.209'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.39%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7517'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.64%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.854'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.151'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.85%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.670'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.16%  '
$ws.Range('D38').Value = '1.237.54'
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.773'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.58%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01795'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.621'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9029'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.001'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '102.53'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '67.09'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.45%  '
$ws.Range('D46').Value = '1.981.80'
$ws.Range('E46').Value = '  -0.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000124'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.76%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5094'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.27%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4098'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('B50').Value = 'XinFinNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07414'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +17.86%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.070'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.30%  '
